$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 774 (pushes existing rows 774:860 down to 775:861)
$ws.Rows(774).Insert()

# Populate the newly inserted row with the new record
$ws.Range("A774").Value = 3
$ws.Range("B774").Value = "Femacal de La Calera"
$ws.Range("C774").Value = "Coquimbo"
$ws.Range("D774").Value = 45212
$ws.Range("E774").Value = 5
$ws.Range("F774").Value = 100112032
$ws.Range("G774").Value = "Zapallo italiano"
$ws.Range("H774").Value = "Sin especificar"
$ws.Range("I774").Value = "Primera"
$ws.Range("J774").Value = 90
$ws.Range("K774").Value = 17000
$ws.Range("L774").Value = 18000
$ws.Range("M774").Value = 17444
$ws.Range("N774").Value = "`$/caja 60 unidades"
$ws.Range("O774").Value = "Región de Arica y Parinacota"
$ws.Range("P774").Value = 291
$ws.Range("Q774").Value = 60
$ws.Range("R774").Value = "Hortaliza"
